# Elimina antiguos EC y agrega nuevos y modifica Antigua BD
#
# Insert a new detail row (period 2509) for worker CC 30871208 "ANAIS TOM
# GUERRERO" right after the existing last detail row (row 80), pushing the
# footer rows (signature line / legal representative name) down by one row,
# and update the summary totals accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 81 (i.e. right after the current last data
# row 80, before the blank rows leading to the footer). This shifts the
# footer rows (85, 86) down to (86, 87). Use a targeted range (not a full
# row) so we don't touch formatting outside the used columns.
$ws.Range("B81:J81").Insert(-4121)  # xlShiftDown

# The new row 81 becomes the new "last" detail row, so it should carry the
# heavier bottom-border formatting that row 80 used to have.
$ws.Range("B80:J80").Copy()
$ws.Range("B81:J81").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 80 is no longer the last row, so it switches to the regular
# detail-row formatting used by the rest of the table (e.g. row 79).
$ws.Range("B79:J79").Copy()
$ws.Range("B80:J80").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new detail row's values.
$ws.Range("B81").Value = "CC"
$ws.Range("C81").Value = "30871208"
$ws.Range("D81").Value = "ANAIS TOM GUERRERO"
$ws.Range("E81").Value = "2509"
$ws.Range("F81").Value = 40000
$ws.Range("G81").Value = 1000000

# Update the summary header: total overdue value (VALOR MORA) and the
# period count (Cant. Periodos) now that one more period was added.
$ws.Range("E11").Value = 3120030
$ws.Range("F13").Value = 51
